$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Nädal 6" sheet by duplicating "Nädal 5" (keeps all the
#    shared layout: header rows, merged cells, styles, formulas, page setup)
#    and appending it as the last tab - same as the author did in Excel via
#    "Move or Copy... > Create a copy".
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("Nädal 5")
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Nädal 6"

# ---------------------------------------------------------------------------
# 2. Fill in week 6's log data on the new sheet.
# ---------------------------------------------------------------------------

# Header date
$new.Range("G3").Value = 43528

# Row 6
$new.Range("B6").Value = 43528
$new.Range("C6").Value = 0.66666666666666663
$new.Range("D6").Value = 0.72916666666666663
$new.Range("G6").Value = "Class "
$new.Range("H6").Value = "Lecture"

# Row 7
$new.Range("B7").Value = 43529
$new.Range("C7").Value = 0.66666666666666663
$new.Range("D7").Value = 0.83333333333333337
$new.Range("E7").Value = 10
$new.Range("G7").Value = "Proge"

# Row 8
$new.Range("B8").Value = 43529
$new.Range("C8").Value = 0.91666666666666663
$new.Range("D8").Value = 1
$new.Range("G8").Value = "Proge"
$new.Range("H8").Value = "Book"

# Row 9
$new.Range("B9").Value = 43530
$new.Range("C9").Value = 0.70833333333333337
$new.Range("D9").Value = 0.875
$new.Range("E9").Value = 20
$new.Range("G9").Value = "Proge"
$new.Range("H9").Value = "Book"

# Row 10
$new.Range("B10").Value = 43531
$new.Range("C10").Value = 0.45833333333333331
$new.Range("D10").Value = 0.5625
$new.Range("E10").Value = 10
$new.Range("G10").Value = "Proge"
$new.Range("H10").Value = "Book"

# Row 11
$new.Range("B11").Value = 43531
$new.Range("C11").Value = 0.58333333333333337
$new.Range("D11").Value = 0.68055555555555547
$new.Range("E11").Value = 20
$new.Range("G11").Value = "Proge"

# Row 12
$new.Range("B12").Value = 43531
$new.Range("C12").Value = 0.6875
$new.Range("D12").Value = 0.82638888888888884
$new.Range("E12").Value = 20
$new.Range("G12").Value = "Proge"
$new.Range("J12").ClearContents()

# Row 13 (new data row - was blank on the template sheet)
$new.Range("A13").Value = 8
$new.Range("B13").Value = 43531
$new.Range("C13").Value = 0.875
$new.Range("D13").Value = 0.97916666666666663
$new.Range("E13").Value = 15
$new.Range("F13").Formula = "=(D13-C13)*24*60 - E13"
$new.Range("G13").Value = "Proge"
$new.Range("H13").Value = "Book"
$new.Range("I13").Value = "x"

# Column B a touch narrower on the new sheet (author hand-resized it).
$new.Columns.Item(2).ColumnWidth = 13.29

# ---------------------------------------------------------------------------
# 3. Selection bookkeeping: the old active sheet ("Nädal 5") keeps a plain
#    selection at F6, the new sheet ("Nädal 6") becomes active/selected with
#    its cursor at K11.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Nädal 5")
$ws5.Activate()
$ws5.Range("F6").Select()

$new.Activate()
$new.Range("K11").Select()
